$d = $word.ActiveDocument

# 1. Replace the placeholder ID text and drop the trailing lone-space run
#    that followed it in the first paragraph.
$d.Content.Find.Execute(
    "**ID__AFFARS_mp_5301_602_2_d_topic_6__ID** ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "**ID__AFFARS_MP5301_602_2_5__ID**", 2)

# 2. Update the first paragraph's formatting: add a paragraph border and
#    change the left indent from 120 to 225 twips.
$p1 = $d.Paragraphs(1)
$p1.Range.ParagraphFormat.LeftIndent = 11.25

$border = $p1.Range.ParagraphFormat.Borders
$border.DistanceFromTop = 5
$border.DistanceFromLeft = 5
$border.DistanceFromBottom = 5
$border.DistanceFromRight = 5
